$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column A duplicates column F (GENE id). Delete column A entirely,
# shifting B:F left into A:E.
$ws.Range("A:A").Delete()

# Fix the shared-string typo: MODEL_CONDITION -> MODELCONDITION
# (after the column shift, this header is now in column D)
$ws.Range("D1").Value = "MODELCONDITION"
